$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 21 de Marzo de 2020 a las 22:46'

$ws.Range("B7").Value = 24786
$ws.Range("C7").Value = 5403
$ws.Range("E7").Value = 24325
$ws.Range("G7").Value = 34
$ws.Range("H7").Value = 290

$ws.Range("B12").Value = 6747
$ws.Range("C12").Value = 1132
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 80

$ws.Range("B15").Value = 2970
$ws.Range("C15").Value = 321
$ws.Range("E15").Value = 2953

$ws.Range("B17").Value = 2157
$ws.Range("C17").Value = 198
$ws.Range("E17").Value = 2144

$ws.Range("B25").Value = 1054
$ws.Range("C25").Value = 47
$ws.Range("E25").Value = 803

$ws.Range("F35").Value = 18

$ws.Range("A107").Value = 'Camerun'
$ws.Range("B107").Value = 40
$ws.Range("C107").Value = 13
$ws.Range("D107").Value = 2
$ws.Range("E107").Value = 38

$ws.Range("A108").Value = 'Liechtenstein'
$ws.Range("B108").Value = 37
$ws.Range("C108").Value = 9
$ws.Range("E108").Value = 37
$ws.Range("F108").Value = 0
$ws.Range("H108").Value = 0

$ws.Range("A109").Value = 'Martinica'
$ws.Range("B109").Value = 32
$ws.Range("D109").Value = 0
$ws.Range("E109").Value = 31
$ws.Range("F109").Value = 7
$ws.Range("H109").Value = 1

$ws.Range("A121").Value = 'Paraguay'
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 0
$ws.Range("F121").Value = 1
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 1

$ws.Range("A122").Value = 'Monaco'
$ws.Range("C122").Value = 7
$ws.Range("D122").Value = 1
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 0

$ws.Range("A140").Value = 'Guinea Ecuatorial'

$ws.Range("A141").Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range("C141").Value = 3

$ws.Range("A142").Value = 'Tanzania'
$ws.Range("C142").Value = 0

$ws.Range("A143").Value = 'Barbados'

$ws.Range("A145").Value = 'Aruba'
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 1
$ws.Range("H145").Value = 0

$ws.Range("A146").Value = 'Guyana'
$ws.Range("D146").Value = 0
$ws.Range("H146").Value = 1

$ws.Range("A147").Value = 'Gabon'
$ws.Range("C147").Value = 1

$ws.Range("A148").Value = 'San Martin (Parte Francesa)'

$ws.Range("A150").Value = 'Bahamas'

$ws.Range("A151").Value = 'Madagascar'
$ws.Range("C151").Value = 0

$ws.Range("A154").Value = 'San Bartolome'

$ws.Range("A158").Value = 'Cabo Verde'
$ws.Range("C158").Value = 2

$ws.Range("A159").Value = 'Republica de Africa Central'

$ws.Range("A163").Value = 'Butan'

$ws.Range("A164").Value = 'Zambia'

$ws.Range("A165").Value = 'Groenlandia'

$ws.Range("A166").Value = 'Fiyi'
$ws.Range("C166").Value = 1

$ws.Range("A167").Value = 'Santa Lucia'

$ws.Range("A168").Value = 'Guinea'

$ws.Range("A169").Value = 'Benin'

$ws.Range("A172").Value = 'Mauritania'

$ws.Range("A173").Value = 'Isla de Man'
$ws.Range("C173").Value = 0

$ws.Range("A176").Value = 'Niger'

$ws.Range("A177").Value = 'Suazilandia'

$ws.Range("A178").Value = 'Montserrat'

$ws.Range("A180").Value = 'Antigua y Barbuda'

$ws.Range("A181").Value = 'San Vicente y las Granadinas'

$ws.Range("A182").Value = 'Santa Sede'

$ws.Range("A183").Value = 'Eritrea'

$ws.Range("A184").Value = 'Papua Nueva Guinea'

$ws.Range("A185").Value = 'Uganda'
$ws.Range("C185").Value = 1

$ws.Range("A186").Value = 'Somalia'
$ws.Range("C186").Value = 0

$ws.Range("A188").Value = 'Republica del Chad'

$ws.Range("A189").Value = 'Timor Oriental'
$ws.Range("C189").Value = 1

$ws.Range("A190").Value = 'San Martin (Parte Holandesa)'
$ws.Range("C190").Value = 0
